$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (borders/font/alignment) from row 13 down to the new rows 14-22
$src = $ws.Range("A13:C13")
$dst = $ws.Range("A14:C22")
$src.Copy($dst)

# Match row heights for the newly added rows to the existing pattern
for ($r = 14; $r -le 22; $r++) {
    $ws.Rows.Item($r).RowHeight = 15.75
}

# Fill in the new application names
$ws.Range("A13").Value = "BBC News"
$ws.Range("A14").Value = "Bubble Shoot"
$ws.Range("A15").Value = "Color Note"
$ws.Range("A16").Value = "Evernote"

# Update the view: scroll so row 4 is at top, and select B17 (as left by the editor)
$ws.Range("B17").Select()
$excel.ActiveWindow.ScrollRow = 4
